$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tValues = @{
    3  = 2023
    4  = 9.6
    5  = 14.6
    6  = 4.7
    7  = 14
    8  = 22.3
    9  = 5.6
    10 = 5.8
    11 = 8.4
    12 = 3.2
    13 = 7.8
    14 = 12.6
    15 = 3
    16 = 20.6
    17 = 28.1
    18 = 13
    19 = 7.9
    20 = 12.2
    21 = 3.6
    22 = 0.7
    23 = 1.4
    25 = 12.2
    26 = 20
    27 = 4.5999999999999996
    28 = 12.9
    29 = 18.899999999999999
    30 = 7.9
    31 = 4.0999999999999996
    32 = 7
    33 = 1.1000000000000001
}

foreach ($row in $tValues.Keys) {
    $ws.Cells.Item($row, 20).Value = $tValues[$row]
}
$ws.Cells.Item(24, 20).Value = "-"

Write-Output "done"
